# feat: add 2022-Q4 data
#
# 1) Duplicate the existing "2022-Q3" sheet (this carries over all of its
#    formatting/structure exactly), drop it in right before "2022-Q3", rename
#    it to "2022-Q4", overwrite its data rows with the new quarter's figures
#    and trim the now-unused trailing rows.
# 2) Insert a new row at the top of the "总计" (summary) sheet for the
#    2022-Q4 totals, shifting the existing quarters down and renumbering the
#    index column.

function Set-TextValue {
    # Forces a value to be stored as text, even when it looks numeric
    # (e.g. a fund code like "002229" or a figure like "11.39" that must
    # keep its original formatting/leading zeros), then strips the
    # leftover "@" number-format style so the cell matches a plain,
    # unstyled text cell.
    param($Cell, $Val)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Val
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: add the "2022-Q4" worksheet before "2022-Q3"
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

$q4Data = @(
    @(0, "002229", "华夏经济转型股票",       "11.39", "85.49", "3.40", "0.3873", 7),
    @(1, "006868", "华夏科技成长股票",       "5.05",  "88.93", "4.12", "0.2081", 2),
    @(2, "013877", "财通资管新能源汽车混合C", "0.74",  "94.54", "4.37", "0.0323", 10),
    @(3, "013876", "财通资管新能源汽车混合A", "0.14",  "94.54", "4.37", "0.0061", 10)
)

$r = 2
foreach ($row in $q4Data) {
    $q4Sheet.Cells.Item($r,1).Value = $row[0]
    Set-TextValue $q4Sheet.Cells.Item($r,2) $row[1]
    Set-TextValue $q4Sheet.Cells.Item($r,3) $row[2]
    Set-TextValue $q4Sheet.Cells.Item($r,4) $row[3]
    Set-TextValue $q4Sheet.Cells.Item($r,5) $row[4]
    Set-TextValue $q4Sheet.Cells.Item($r,6) $row[5]
    Set-TextValue $q4Sheet.Cells.Item($r,7) $row[6]
    $q4Sheet.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# the copied "2022-Q3" sheet has 9 data rows (rows 2-10); 2022-Q4 only has 4
# (rows 2-5), so drop the now-stale trailing rows.
$q4Sheet.Range("A6:H10").Delete()

# ---------------------------------------------------------------------------
# Step 2: update the "总计" sheet with the new 2022-Q4 row
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2").EntireRow.Insert()

# carry the existing formatting (index-column style) down into the new row
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q4"
$totalSheet.Cells.Item(2,3).Value = 4
$totalSheet.Cells.Item(2,4).Value = 0.63

# renumber the index column (A) for the rows that shifted down
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(6,1).Value = 4
